# Update cryptos list with latest price/volume data (GitHub Actions sync)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '51.083.37'
$ws.Range("E2").Value = '  +0.20%  '

# Row 3
$ws.Range("D3").Value = '2.956.54'
$ws.Range("E3").Value = '  +0.41%  '

# Row 4
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").Value = '''379.87'
$ws.Range("E5").Value = '  +1.08%  '

# Row 6
$ws.Range("D6").Value = '''102.04'
$ws.Range("E6").Value = '  +0.31%  '

# Row 7
$ws.Range("D7").Value = '''0.544'
$ws.Range("E7").Value = '  +1.59%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("E9").Value = '  +0.33%  '

# Row 10
$ws.Range("D10").Value = '''36.36'
$ws.Range("E10").Value = '  +0.11%  '

# Row 11
$ws.Range("E11").Value = '  -1.29%  '

# Row 12
$ws.Range("D12").Value = '''0.0855'
$ws.Range("E12").Value = '  +1.84%  '

# Row 13
$ws.Range("D13").Value = '3.418.71'
$ws.Range("E13").Value = '  +0.40%  '

# Row 14
$ws.Range("D14").Value = '''7.79'
$ws.Range("E14").Value = '  +4.58%  '

# Row 15
$ws.Range("D15").Value = '''18.28'
$ws.Range("E15").Value = '  +1.89%  '

# Row 16
$ws.Range("D16").Value = '2.950.98'
$ws.Range("E16").Value = '  +0.49%  '

# Row 17
$ws.Range("D17").Value = '''11.14'
$ws.Range("E17").Value = '  +3.06%  '

# Row 18
$ws.Range("E18").Value = '  +1.34%  '

# Row 19
$ws.Range("D19").Value = '51.141.07'
$ws.Range("E19").Value = '  +0.52%  '

# Row 20
$ws.Range("D20").Value = '''3.14'
$ws.Range("E20").Value = '  -0.42%  '

# Row 21
$ws.Range("E21").Value = '  -2.34%  '

# Row 22
$ws.Range("D22").Value = '0.0₃0959'
$ws.Range("E22").Value = '  +0.55%  '

# Row 23
$ws.Range("D23").Value = '''70.39'
$ws.Range("E23").Value = '  +2.79%  '

# Row 24
$ws.Range("D24").Value = '''3.28'
$ws.Range("E24").Value = '  +4.57%  '

# Row 25
$ws.Range("D25").Value = '''266.95'
$ws.Range("E25").Value = '  +0.99%  '

# Row 26
$ws.Range("D26").Value = '''7.82'
$ws.Range("E26").Value = '  -5.13%  '

# Row 27
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = '''7.18'
$ws.Range("E27").Value = '  -8.44%  '

# Row 28
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = '''0.999'
$ws.Range("E28").Value = '  -0.09%  '

# Row 29
$ws.Range("D29").Value = '''25.82'
$ws.Range("E29").Value = '  +0.68%  '

# Row 30
$ws.Range("D30").Value = '''0.164'
$ws.Range("E30").Value = '  -3.29%  '

# Row 31
$ws.Range("E31").Value = '  -0.52%  '

# Row 32
$ws.Range("D32").Value = '''10.27'
$ws.Range("E32").Value = '  +3.88%  '

# Row 33
$ws.Range("D33").Value = '''51.06'
$ws.Range("E33").Value = '  +0.68%  '

# Row 34
$ws.Range("D34").Value = '''34.24'
$ws.Range("E34").Value = '  +1.97%  '

# Row 35
$ws.Range("E35").Value = '  +1.90%  '

# Row 36
$ws.Range("D36").Value = '''0.0434'
$ws.Range("E36").Value = '  -1.80%  '

# Row 37
$ws.Range("E37").Value = '  +0.07%  '

# Row 38
$ws.Range("E38").Value = '  +4.57%  '

# Row 39
$ws.Range("E39").Value = '  +0.40%  '

# Row 40
$ws.Range("D40").Value = '''1.83'
$ws.Range("E40").Value = '  +2.13%  '

# Row 41
$ws.Range("E41").Value = '  +0.19%  '

# Row 42
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").Value = '''124.66'
$ws.Range("E42").Value = '  +3.46%  '

# Row 43
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '''2.50'
$ws.Range("E43").Value = '  -0.92%  '

# Row 44
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").Value = '''3.52'
$ws.Range("E44").Value = '  +6.61%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''21.38'
$ws.Range("E45").Value = '  +1.22%  '

# Row 46
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '''2.02'
$ws.Range("E46").Value = '  -0.15%  '

# Row 47
$ws.Range("B47").Value = 'TheGraph'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D47").Value = '''0.271'
$ws.Range("E47").Value = '  -5.03%  '

# Row 48
$ws.Range("D48").Value = '''2.36'
$ws.Range("E48").Value = '  +2.75%  '

# Row 49
$ws.Range("D49").Value = '2.038.46'
$ws.Range("E49").Value = '  +2.20%  '

# Row 50
$ws.Range("D50").Value = '''0.0321'
$ws.Range("E50").Value = '  -5.64%  '

# Row 51
$ws.Range("E51").Value = '  +6.27%  '

# Reset style on text-forced numeric-looking price cells so the
# quote-prefix marker doesn't linger as cell formatting metadata
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
